$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 270 so the two rows that used
# to be 270/271 become 272/273, making room for two "new" rows at 270/271.
$ws.Rows.Item(270).Insert()
$ws.Rows.Item(270).Insert()

# Row 267 (date + price band changed)
$ws.Range("D267").Value = 44448
$ws.Range("J267").Value = 440
$ws.Range("K267").Value = 1500
$ws.Range("L267").Value = 1500
$ws.Range("M267").Value = 1500
$ws.Range("P267").Value = 500

# Row 268 (date + price band changed)
$ws.Range("D268").Value = 44448
$ws.Range("J268").Value = 350
$ws.Range("K268").Value = 1000
$ws.Range("L268").Value = 1000
$ws.Range("M268").Value = 1000
$ws.Range("P268").Value = 333

# Row 269 (date rolled back, values restored to the old row 267 figures)
$ws.Range("D269").Value = 44167
$ws.Range("J269").Value = 610
$ws.Range("K269").Value = 800
$ws.Range("L269").Value = 1000
$ws.Range("M269").Value = 885
$ws.Range("P269").Value = 295

# Row 270 (date rolled back, quality/origin/values restored to the old row 268 figures)
$ws.Range("A270").Value = 6
$ws.Range("B270").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C270").Value = "Metropolitana"
$ws.Range("D270").Value = 44167
$ws.Range("E270").Value = 13
$ws.Range("F270").Value = 100112039
$ws.Range("G270").Value = "Ciboulette"
$ws.Range("H270").Value = "Sin especificar"
$ws.Range("I270").Value = "Segunda"
$ws.Range("J270").Value = 320
$ws.Range("K270").Value = 700
$ws.Range("L270").Value = 800
$ws.Range("M270").Value = 744
$ws.Range("N270").Value = "$/docena de atados"
$ws.Range("O270").Value = "Región Metropolitana"
$ws.Range("P270").Value = 248
$ws.Range("Q270").Value = 3
$ws.Range("R270").Value = "Hortaliza"

# Row 271 (new row, values matching the old row 269 figures)
$ws.Range("A271").Value = 6
$ws.Range("B271").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C271").Value = "Metropolitana"
$ws.Range("D271").Value = 44238
$ws.Range("E271").Value = 13
$ws.Range("F271").Value = 100112039
$ws.Range("G271").Value = "Ciboulette"
$ws.Range("H271").Value = "Sin especificar"
$ws.Range("I271").Value = "Primera"
$ws.Range("J271").Value = 630
$ws.Range("K271").Value = 700
$ws.Range("L271").Value = 800
$ws.Range("M271").Value = 760
$ws.Range("N271").Value = "$/docena de atados"
$ws.Range("O271").Value = "Región Metropolitana"
$ws.Range("P271").Value = 253
$ws.Range("Q271").Value = 3
$ws.Range("R271").Value = "Hortaliza"

# Row 272 (old row 270, now shifted down two positions, values unchanged)
$ws.Range("A272").Value = 6
$ws.Range("B272").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C272").Value = "Metropolitana"
$ws.Range("D272").Value = 44399
$ws.Range("E272").Value = 13
$ws.Range("F272").Value = 100112039
$ws.Range("G272").Value = "Ciboulette"
$ws.Range("H272").Value = "Sin especificar"
$ws.Range("I272").Value = "Primera"
$ws.Range("J272").Value = 700
$ws.Range("K272").Value = 1800
$ws.Range("L272").Value = 2000
$ws.Range("M272").Value = 1886
$ws.Range("N272").Value = "$/docena de atados"
$ws.Range("O272").Value = "Provincia de Chacabuco"
$ws.Range("P272").Value = 629
$ws.Range("Q272").Value = 3
$ws.Range("R272").Value = "Hortaliza"

# Row 273 (old row 271, now shifted down two positions, values unchanged)
$ws.Range("A273").Value = 6
$ws.Range("B273").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C273").Value = "Metropolitana"
$ws.Range("D273").Value = 44400
$ws.Range("E273").Value = 13
$ws.Range("F273").Value = 100112039
$ws.Range("G273").Value = "Ciboulette"
$ws.Range("H273").Value = "Sin especificar"
$ws.Range("I273").Value = "Primera"
$ws.Range("J273").Value = 700
$ws.Range("K273").Value = 1800
$ws.Range("L273").Value = 2000
$ws.Range("M273").Value = 1886
$ws.Range("N273").Value = "$/docena de atados"
$ws.Range("O273").Value = "Provincia de Chacabuco"
$ws.Range("P273").Value = 629
$ws.Range("Q273").Value = 3
$ws.Range("R273").Value = "Hortaliza"
